$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '43.799.18'
$ws.Range('E2').Value = '  -0.18%  '
$ws.Range('D3').Value = '2.337.24'
$ws.Range('E3').Value = '  -0.76%  '
$ws.Range('E4').Value = '  +0.05%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '238.72'
$ws.Range('E5').Value = '  -1.44%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.657'
$ws.Range('E6').Value = '  -5.15%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '71.81'
$ws.Range('E7').Value = '  -7.03%  '
$ws.Range('E8').Value = '  -0.04%  '
$ws.Range('E9').Value = '  -5.10%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.0988'
$ws.Range('E10').Value = '  -3.52%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '57.89'
$ws.Range('E11').Value = '  +0.88%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '32.12'
$ws.Range('E12').Value = '  -5.26%  '
$ws.Range('E13').Value = '  -0.61%  '
$ws.Range('E14').Value = '  -6.15%  '
$ws.Range('D15').Value = '2.685.40'
$ws.Range('E15').Value = '  -0.66%  '
$ws.Range('E16').Value = '  -5.67%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.896'
$ws.Range('E17').Value = '  -3.36%  '
$ws.Range('D18').Value = '2.340.12'
$ws.Range('E18').Value = '  -0.56%  '
$ws.Range('D19').Value = '43.669.22'
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '0.0000100'
$ws.Range('E20').Value = '  -2.88%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '77.70'
$ws.Range('E21').Value = '  +0.11%  '
$ws.Range('E22').Value = '  -1.34%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '250.78'
$ws.Range('E23').Value = '  -2.21%  '
$ws.Range('B24').Value = 'Dai'
$ws.Range('C24').Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '1.00'
$ws.Range('E24').Value = '  +0.03%  '
$ws.Range('B25').Value = 'ImmutableX'
$ws.Range('C25').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '1.90'
$ws.Range('E25').Value = '  +5.60%  '
$ws.Range('E26').Value = '  +2.76%  '
$ws.Range('E27').Value = '  -2.37%  '
$ws.Range('E28').Value = '  -6.79%  '
$ws.Range('E29').Value = '  -1.02%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '176.28'
$ws.Range('E30').Value = '  +0.78%  '
$ws.Range('E31').Value = '  -4.15%  '
$ws.Range('E32').Value = '  -2.48%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.133'
$ws.Range('E33').Value = '  -2.10%  '
$ws.Range('E34').Value = '  -3.33%  '
$ws.Range('E35').Value = '  -5.32%  '
$ws.Range('E36').Value = '  -1.92%  '
$ws.Range('E37').Value = '  -2.38%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '5.93'
$ws.Range('E38').Value = '  +31.70%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '6.36'
$ws.Range('E39').Value = '  -1.83%  '
$ws.Range('E40').Value = '  -3.49%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.0269'
$ws.Range('E41').Value = '  -3.13%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '66.40'
$ws.Range('E42').Value = '  +19.82%  '
$ws.Range('B43').Value = 'Cronos'
$ws.Range('C43').Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.107'
$ws.Range('E43').Value = '  +4.15%  '
$ws.Range('B44').Value = 'FraxShare'
$ws.Range('C44').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '9.11'
$ws.Range('E44').Value = '  +0.50%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '18.64'
$ws.Range('E45').Value = '  -4.63%  '
$ws.Range('E46').Value = '  -3.82%  '
$ws.Range('E47').Value = '  +0.09%  '
$ws.Range('E48').Value = '  -4.30%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '2.41'
$ws.Range('E49').Value = '  -5.20%  '
$ws.Range('E50').Value = '  +3.31%  '
$ws.Range('E51').Value = '  -4.12%  '
